$wb = $excel.ActiveWorkbook

# Update the "repaymentstrategy" value on the ProductLoanInput sheet (row 17, col B)
# from "RBI (India)" to "Overdue/Due Fee/Int,Principal"
$ws = $wb.Worksheets.Item("ProductLoanInput")
$ws.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

# Move/record the active selection to B17 (single cell) as captured in the edit
$ws.Range("B17").Select()
